$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.683.15"
$ws.Range("E2").Value = "  -0.36%  "
$ws.Range("D3").Value = "2.370.96"
$ws.Range("E3").Value = "  -3.67%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "541.82"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.11%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "140.08"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.79%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.543"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -9.47%  "
$ws.Range("D9").Value = "2.365.62"
$ws.Range("E9").Value = "  -3.86%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.104"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.53%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.154"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.36%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.31"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.05%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.342"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.75%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "25.33"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.63%  "
$ws.Range("D15").Value = "2.803.63"
$ws.Range("E15").Value = "  -3.59%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000161"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.16%  "
$ws.Range("D17").Value = "60.614.41"
$ws.Range("E17").Value = "  -0.34%  "
$ws.Range("D18").Value = "2.374.41"
$ws.Range("E18").Value = "  -3.59%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.56"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -4.34%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.08"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.23%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "315.15"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.15%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.66"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.73%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.999"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.02%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.80"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.54%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "62.74"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.01%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.19%  "
$ws.Range("D27").Value = "2.495.73"
$ws.Range("E27").Value = "  -3.35%  "
$ws.Range("B28").Value = "PEPE"
$ws.Range("C28").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D28").Value = "0.0₃0916"
$ws.Range("E28").Value = "  -6.33%  "
$ws.Range("B29").Value = "Aptos"
$ws.Range("C29").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.63"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.35%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "513.57"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -3.88%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.41"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -4.48%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.91"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -4.55%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.144"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -4.50%  "
$ws.Range("E34").Value = "  -3.55%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.54"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.81%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.00"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.05%  "
$ws.Range("B37").Value = "NEARProtocol"
$ws.Range("C37").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.61"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -4.79%  "
$ws.Range("B38").Value = "RenderToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.41"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -7.69%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.372"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.85%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "17.97"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.05%  "
$ws.Range("B41").Value = "USDe"
$ws.Range("C41").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.00"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.16%  "
$ws.Range("B42").Value = "Stacks"
$ws.Range("C42").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.70"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.31%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "136.43"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -6.42%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "40.23"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.92%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.21"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "138.20"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -5.92%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.52"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.14%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "20.10"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.90%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0514"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.97%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.573"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.85%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0906"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.36%  "
